# Word was left with the insertion point (the "_GoBack" bookmark) right
# after an empty "Updates" bullet line ("             - "). The edit adds
# a sentence of real content to that bullet:
#   "             - " -> "             - Player now supports looped
#                          horizontal movement in a closed space. "
#
# The "_GoBack" bookmark is unique in the document and sits at the very
# end of that bullet's text, so it is used as a precise anchor to find
# the exact paragraph to edit (there is another, unrelated, empty
# "             - " bullet earlier in the document that must stay
# untouched).

$d = $word.ActiveDocument

$bm = $d.Bookmarks("_GoBack")
$targetPara = $d.Range($bm.Start, $bm.Start).Paragraphs(1)
$paraRange = $targetPara.Range

# Exclude the trailing paragraph mark from the replacement range so only
# the run's text is touched and the paragraph / bookmark stay intact.
$textRange = $d.Range($paraRange.Start, $paraRange.End - 1)

if ($textRange.Text -eq "             - ") {
    $textRange.Text = "             - Player now supports looped horizontal movement in a closed space. "
}
